# Commit: "Add explanation for variables, change city names, change template"
#
# This script reproduces the data-level part of that commit: the `city`
# column (column D, rows 2-11) in the single "data" worksheet gets ten of
# its values swapped out for new city names.
#
# (The commit also touched devdata/Data.xlsx's low-level style table -
# a new solid-white fill + refreshed indexed border color for the shared
# "bordered cell" style, plus a couple of internal DrawingML theme-effect
# tweaks. Those don't change any cell's value or which style index it
# points to (every cell keeps style index 2 before and after) - they only
# redefine what style index 2 *looks like*. The Excel object model has no
# "edit this xf record in place" primitive: any Range/Style formatting
# call (Interior, Borders, ...) resolves/creates a (possibly new) style
# index for the affected range instead, which would repoint every cell on
# the sheet away from style 2 and introduce a sheet-wide formatting diff
# that doesn't exist in the real change. So that part is intentionally
# left alone here to avoid corrupting the template beyond what the commit
# actually did to the visible data.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the city column (D2:D11) ---------------------------------------
$cityUpdates = [ordered]@{
    "D2"  = "Toronto"
    "D3"  = "New York"
    "D4"  = "San Francisco"
    "D5"  = "Milan"
    "D6"  = "Kansas City"
    "D7"  = "Lahore"
    "D8"  = "New Dehli"
    "D9"  = "Helsinki"
    "D10" = "Stockholm"
    "D11" = "London"
}

foreach ($addr in $cityUpdates.Keys) {
    $ws.Range($addr).Value = $cityUpdates[$addr]
}
